
# Update the "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
# Mapping of row -> new value is identical for both sheets (sheet "全部类型" simply
# has one extra inserted row, shifting everything below it by one).

$wb = $excel.ActiveWorkbook

$updates = @{
    3229 = 3237
    734  = 736
    122  = 124
    6886 = 6890
    2091 = 2112
    28   = 29
    76   = 77
    23   = 24
    72   = 73
    151  = 153
    189  = 192
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $dim = $ws.UsedRange
    $lastRow = $dim.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 6)  # column F
        $val = $cell.Value2
        if ($updates.ContainsKey($val)) {
            $cell.Value2 = $updates[$val]
        }
    }
}
